$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.85780930519104
$ws.Range("B1").Value = 0.7002959847450256
$ws.Range("C1").Value = 3.366165161132812
$ws.Range("D1").Value = 3.017059564590454
$ws.Range("E1").Value = 1.339669227600098
